$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3605.3572
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 4247.5
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 12742.5
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -13282.5
$ws.Range("H73").Value = 3605.3572
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 4247.5
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 12742.5
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -14614.5
$ws.Range("H113").Value = 4872.6113
$ws.Range("I113").Value = 2787
$ws.Range("K113").Value = 2787
$ws.Range("M113").Value = 467
$ws.Range("H125").Value = 1733
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1733
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 15597
$ws.Range("N125").Value = -20517
$ws.Range("M125").ClearContents()
$ws.Range("H129").Value = 883.9299999999999
$ws.Range("J129").Value = 912.8617
$ws.Range("L129").Value = 2738.5851
$ws.Range("N129").Value = -12738.5851

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5964.9653
$ws.Range("I32").Value = 5022.595
$ws.Range("K32").Value = 5022.595
$ws.Range("M32").Value = -4735.595
$ws.Range("H63").Value = 9237609
$ws.Range("I63").Value = 13853574
$ws.Range("J63").Value = 5680
$ws.Range("K63").Value = 13853574
$ws.Range("L63").Value = 5680
$ws.Range("M63").Value = -13852888
$ws.Range("N63").Value = -7052
$ws.Range("H66").Value = 9237609
$ws.Range("I66").Value = 13853574
$ws.Range("J66").Value = 5680
$ws.Range("K66").Value = 69267870
$ws.Range("L66").Value = 28400
$ws.Range("M66").Value = -69264438
$ws.Range("N66").Value = -35264
$ws.Range("H74").Value = 10442.4
$ws.Range("I74").Value = 11660.571
$ws.Range("J74").Value = 7600
$ws.Range("K74").Value = 11660.571
$ws.Range("L74").Value = 7600
$ws.Range("M74").Value = -10786.571
$ws.Range("N74").Value = -9348
$ws.Range("H77").Value = 10442.4
$ws.Range("I77").Value = 11660.571
$ws.Range("J77").Value = 7600
$ws.Range("K77").Value = 58302.855
$ws.Range("L77").Value = 38000
$ws.Range("M77").Value = -53934.855
$ws.Range("N77").Value = -46736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 21278.125
$ws.Range("J81").Value = 21278.125
$ws.Range("L81").Value = 21278.125
$ws.Range("N81").Value = -23400.125
$ws.Range("H84").Value = 21278.125
$ws.Range("J84").Value = 21278.125
$ws.Range("L84").Value = 63834.375
$ws.Range("N84").Value = -74442.375
$ws.Range("H86").Value = 1879.7715
$ws.Range("I86").Value = 1648.6
$ws.Range("J86").Value = 2457.7
$ws.Range("K86").Value = 1648.6
$ws.Range("L86").Value = 2457.7
$ws.Range("M86").Value = -525.5999999999999
$ws.Range("N86").Value = -4703.7
$ws.Range("H89").Value = 1879.7715
$ws.Range("I89").Value = 1648.6
$ws.Range("J89").Value = 2457.7
$ws.Range("K89").Value = 8243
$ws.Range("L89").Value = 12288.5
$ws.Range("M89").Value = -2627
$ws.Range("N89").Value = -23520.5
$ws.Range("H107").Value = 995.8570999999999
$ws.Range("I107").Value = 1017.6316
$ws.Range("J107").Value = 789
$ws.Range("K107").Value = 1017.6316
$ws.Range("L107").Value = 789
$ws.Range("M107").Value = 902.3684
$ws.Range("N107").Value = -4629
$ws.Range("H134").Value = 1969.119
$ws.Range("I134").Value = 1479.5641
$ws.Range("J134").Value = 8333.333000000001
$ws.Range("K134").Value = 4438.692300000001
$ws.Range("L134").Value = 24999.999
$ws.Range("M134").Value = -1903.692300000001
$ws.Range("N134").Value = -30069.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 814.53845
$ws.Range("I22").Value = 262.375
$ws.Range("J22").Value = 1698
$ws.Range("K22").Value = 262.375
$ws.Range("L22").Value = 1698
$ws.Range("M22").Value = 87.625
$ws.Range("N22").Value = -2398
$ws.Range("H62").Value = 5999.2
$ws.Range("I62").Value = 2995
$ws.Range("K62").Value = 2995
$ws.Range("M62").Value = -2371
$ws.Range("H65").Value = 5999.2
$ws.Range("I65").Value = 2995
$ws.Range("K65").Value = 14975
$ws.Range("M65").Value = -11855
$ws.Range("H86").Value = 7244.75
$ws.Range("I86").Value = 6326.6665
$ws.Range("J86").Value = 9999
$ws.Range("K86").Value = 6326.6665
$ws.Range("L86").Value = 9999
$ws.Range("M86").Value = -5203.6665
$ws.Range("N86").Value = -12245
$ws.Range("H89").Value = 7244.75
$ws.Range("I89").Value = 6326.6665
$ws.Range("J89").Value = 9999
$ws.Range("K89").Value = 31633.3325
$ws.Range("L89").Value = 49995
$ws.Range("M89").Value = -26017.3325
$ws.Range("N89").Value = -61227
$ws.Range("H99").Value = 7411002.5
$ws.Range("I99").Value = 14287540
$ws.Range("J99").Value = 5500.231
$ws.Range("K99").Value = 14287540
$ws.Range("L99").Value = 5500.231
$ws.Range("M99").Value = -14286042
$ws.Range("N99").Value = -8496.231
$ws.Range("H105").Value = 3202
$ws.Range("I105").Value = 2670
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 2670
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -923
$ws.Range("N105").Value = -7494
$ws.Range("H107").Value = 589.5238000000001
$ws.Range("I107").Value = 515.6667
$ws.Range("J107").Value = 1032.6666
$ws.Range("K107").Value = 515.6667
$ws.Range("L107").Value = 1032.6666
$ws.Range("M107").Value = 1404.3333
$ws.Range("N107").Value = -4872.6666
$ws.Range("H126").Value = 7411002.5
$ws.Range("I126").Value = 14287540
$ws.Range("J126").Value = 5500.231
$ws.Range("K126").Value = 42862620
$ws.Range("L126").Value = 16500.693
$ws.Range("M126").Value = -42860150
$ws.Range("N126").Value = -21440.693
$ws.Range("H132").Value = 1701.6842
$ws.Range("I132").Value = 1176.375
$ws.Range("K132").Value = 3529.125
$ws.Range("M132").Value = -999.125
$ws.Range("H134").Value = 3530.25
$ws.Range("I134").Value = 3565.8438
$ws.Range("J134").Value = 3435.3333
$ws.Range("K134").Value = 10697.5314
$ws.Range("L134").Value = 10305.9999
$ws.Range("M134").Value = -8162.5314
$ws.Range("N134").Value = -15375.9999
$ws.Range("H141").Value = 30119.2
$ws.Range("J141").Value = 30119.2
$ws.Range("L141").Value = 30119.2
$ws.Range("N141").Value = -40479.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3215.4
$ws.Range("I137").Value = 997.25
$ws.Range("J137").Value = 4694.1665
$ws.Range("K137").Value = 2991.75
$ws.Range("L137").Value = 14082.4995
$ws.Range("M137").Value = 2108.25
$ws.Range("N137").Value = -24282.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1980.579
$ws.Range("I113").Value = 2315.4546
$ws.Range("K113").Value = 2315.4546
$ws.Range("M113").Value = -145.4546
$ws.Range("H122").Value = 2849.3333
$ws.Range("I122").Value = 2296
$ws.Range("J122").Value = 6169.3335
$ws.Range("K122").Value = 6888
$ws.Range("L122").Value = 18508.0005
$ws.Range("M122").Value = -4438
$ws.Range("N122").Value = -23408.0005
$ws.Range("H140").Value = 38518.57
$ws.Range("J140").Value = 38518.57
$ws.Range("L140").Value = 38518.57
$ws.Range("N140").Value = -48878.57

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 49331.24
$ws.Range("I22").Value = 72618
$ws.Range("J22").Value = 2757.7144
$ws.Range("K22").Value = 72618
$ws.Range("L22").Value = 2757.7144
$ws.Range("M22").Value = -72323
$ws.Range("N22").Value = -3347.7144
$ws.Range("H27").Value = 49331.24
$ws.Range("I27").Value = 72618
$ws.Range("J27").Value = 2757.7144
$ws.Range("K27").Value = 72618
$ws.Range("L27").Value = 2757.7144
$ws.Range("M27").Value = -72511
$ws.Range("N27").Value = -2971.7144
$ws.Range("H40").Value = 6621.1113
$ws.Range("I40").Value = 5124.3335
$ws.Range("J40").Value = 8492.083000000001
$ws.Range("K40").Value = 5124.3335
$ws.Range("L40").Value = 8492.083000000001
$ws.Range("M40").Value = -4988.3335
$ws.Range("N40").Value = -8764.083000000001
$ws.Range("H55").Value = 411
$ws.Range("I55").Value = 290.6
$ws.Range("J55").Value = 511.33334
$ws.Range("K55").Value = 290.6
$ws.Range("L55").Value = 511.33334
$ws.Range("M55").Value = -117.6
$ws.Range("N55").Value = -857.33334
$ws.Range("H122").Value = 2752.182
$ws.Range("I122").Value = 1752.625
$ws.Range("K122").Value = 5257.875
$ws.Range("M122").Value = -2807.875
$ws.Range("H132").Value = 8245.522000000001
$ws.Range("I132").Value = 7415.8687
$ws.Range("J132").Value = 13500
$ws.Range("K132").Value = 22247.6061
$ws.Range("L132").Value = 40500
$ws.Range("M132").Value = -19717.6061
$ws.Range("N132").Value = -45560
$ws.Range("H134").Value = 39986.332
$ws.Range("J134").Value = 39986.332
$ws.Range("L134").Value = 39986.332
$ws.Range("N134").Value = -50126.332
$ws.Range("H136").Value = 3217.1936
$ws.Range("I136").Value = 1407.15
$ws.Range("J136").Value = 6508.1816
$ws.Range("K136").Value = 4221.450000000001
$ws.Range("L136").Value = 19524.5448
$ws.Range("M136").Value = -1671.450000000001
$ws.Range("N136").Value = -24624.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 37666.668
$ws.Range("I62").Value = 6500
$ws.Range("J62").Value = 100000
$ws.Range("K62").Value = 6500
$ws.Range("L62").Value = 100000
$ws.Range("M62").Value = -5876
$ws.Range("N62").Value = -101248
$ws.Range("H65").Value = 37666.668
$ws.Range("I65").Value = 6500
$ws.Range("J65").Value = 100000
$ws.Range("K65").Value = 32500
$ws.Range("L65").Value = 500000
$ws.Range("M65").Value = -29380
$ws.Range("N65").Value = -506240
$ws.Range("H126").Value = 1809.7097
$ws.Range("I126").Value = 1212.8182
$ws.Range("J126").Value = 3268.7778
$ws.Range("K126").Value = 3638.4546
$ws.Range("L126").Value = 9806.3334
$ws.Range("M126").Value = -1168.4546
$ws.Range("N126").Value = -14746.3334
